$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new column at V, shifting the existing V:Y ("style list",
# "media::image::language", "media::video::language",
# "media::audio::language") one column to the right (W:Z).
$ws.Columns("V").Insert()

# Add the new "server_calculation" header in the freshly inserted V1,
# matching the bold header style used by the rest of row 1.
$ws.Range("V1").Value = "server_calculation"
$ws.Range("V1").Font.Bold = $true

# Update the view: scroll the frozen pane so the new column is visible,
# then select the new active cell (mirrors what the author's Excel session
# recorded: topLeftCell moved from C2 to T2, active cell from A2 to V2).
$excel.ActiveWindow.ScrollColumn = 20
$ws.Range("V2").Select()
